# Apply updated cryptocurrency price/volume data to sheet1
# (commit: "Updated cryptos list on Thu Oct 19 03:20:19 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.282.24"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.551.71"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.17"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.42"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.771.88"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "1.551.10"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "28.280.06"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.508"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.44"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.33"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.30"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "0.0₃0674"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.81"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  -4.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.89"
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.79"
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.21"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0467"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.06"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").Value = "1.384.77"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.57"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.511"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.775"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0465"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.42"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.79"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.907"
$ws.Range("E47").Value = "  -6.58%  "
$ws.Range("D48").Value = "1.686.08"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.38"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "41.72"
$ws.Range("E51").Value = "  +5.27%  "
